# Applies attendance-status updates to Sheet1 of the active workbook.
# For each student-date row, mark the appropriate status column with 1
# (it was previously left at 0), matching the updated attendance data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> cell(s) that flip from 0 to 1 (column letter => value)
$updates = @{
    3  = @("G3", "H3")
    4  = @("D4", "E4")
    5  = @("D5", "E5")
    6  = @("H6")
    7  = @("H7")
    8  = @("H8")
    9  = @("H9")
    10 = @("G10", "H10")
    11 = @("D11", "E11")
    12 = @("D12", "E12")
    13 = @("H13")
    14 = @("H14")
    15 = @("H15")
    16 = @("H16")
    17 = @("H17")
    18 = @("H18")
}

foreach ($row in $updates.Keys) {
    foreach ($cellRef in $updates[$row]) {
        $ws.Range($cellRef).Value = 1
    }
}
